$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: After paragraph "DbgMe and DbgMeAsStr in Model" (paragraph 10),
# insert two new list paragraphs:
#   "Main parameters, testability, maintenance" (ilvl 0)
#   "Test Mode in Main "                        (ilvl 1)
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item(11)
$newPara1.Range.ListFormat.ListLevelNumber = 1
$newPara1.Range.Text = "Main parameters, testability, maintenance"

$newPara1 = $d.Paragraphs.Item(11)
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(12)
$newPara2.Range.ListFormat.ListLevelNumber = 2
$newPara2.Range.Text = "Test Mode in Main "

Write-Host "=== paragraphs after step 1 ==="
for ($i = 1; $i -le 17; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host $i ": [" $p.Range.Text "]"
}

# ---------------------------------------------------------------------------
# Step 2: the (originally highlighted) paragraph 13 "Main parameters,
# testability, maintenance" becomes "Test"
# ---------------------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Find.Execute("Main parameters, testability, maintenance", $true, $false, $false, $false, $false, $true, 1, $false, "Test", 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: the (originally highlighted) paragraph 14 "Test Mode in Main "
# becomes "Unit test"; a new highlighted paragraph "assert" (ilvl 1) follows
# ---------------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Find.Execute("Test Mode in Main ", $true, $false, $false, $false, $false, $true, 1, $false, "Unit test", 2) | Out-Null

$p14 = $d.Paragraphs.Item(14)
$p14.Range.InsertParagraphAfter()
$newPara3 = $d.Paragraphs.Item(15)
$newPara3.Range.Font.HighlightColorIndex = 7
$newPara3.Range.Text = "assert"

Write-Host "=== paragraphs after step 2+3 ==="
for ($i = 1; $i -le 19; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host $i ": [" $p.Range.Text "]"
}

# ---------------------------------------------------------------------------
# Step 4: remove one of the two consecutive empty paragraphs
# ---------------------------------------------------------------------------
$pEmpty = $d.Paragraphs.Item(16)
Write-Host "pEmpty to delete: [" $pEmpty.Range.Text "]"
$pEmpty.Range.Delete()

Write-Host "=== paragraphs after step 4 ==="
for ($i = 1; $i -le 19; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host $i ": [" $p.Range.Text "]"
}
